$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

$totalWs = $wb.Worksheets.Item(1)     # "总计"
$q2Ws    = $wb.Worksheets.Item(2)     # "2022-Q2" (existing data sheet)

# ---------------------------------------------------------------------------
# 1. Update the "总计" (totals) sheet: change the existing Q2 row into the new
#    Q4 row, then re-append a fresh row underneath with the original Q2 data.
# ---------------------------------------------------------------------------
$totalWs.Range("B2").Value = "2022-Q4"
$totalWs.Range("D2").Value = 0.27

$totalWs.Range("A3").Value = 1
$totalWs.Range("B3").Value = "2022-Q2"
$totalWs.Range("C3").Value = 1
$totalWs.Range("D3").Value = 0.17

# Give A3 the same style as A2 ("总计" header/index style).
$totalWs.Range("A2").Copy()
$totalWs.Range("A3").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------------
# 2. Insert a brand-new worksheet for the "2022-Q4" fund holdings, positioned
#    right before the existing "2022-Q2" sheet (so the tab order becomes
#    总计, 2022-Q4, 2022-Q2).
# ---------------------------------------------------------------------------
$q4Ws = $wb.Worksheets.Add($q2Ws)
$q4Ws.Name = "2022-Q4"

# Header row
$q4Ws.Range("B1").Value = "基金代码"
$q4Ws.Range("C1").Value = "基金名称"
$q4Ws.Range("D1").Value = "基金规模"
$q4Ws.Range("E1").Value = "股票总仓位"
$q4Ws.Range("F1").Value = "仓位占比"
$q4Ws.Range("G1").Value = "持有市值(亿元)"
$q4Ws.Range("H1").Value = "仓位排名"

# Data row
$q4Ws.Range("A2").Value = 0
$q4Ws.Range("B2").Value = "'007207"
$q4Ws.Range("C2").Value = "华夏常阳三年定期开放混合"
$q4Ws.Range("D2").Value = "'6.99"
$q4Ws.Range("E2").Value = "'73.26"
$q4Ws.Range("F2").Value = "'3.82"
$q4Ws.Range("G2").Value = "'0.2670"
$q4Ws.Range("H2").Value = 3

# Match the "总计" sheet's header/index style (already style index 2 in the
# template) on the header row and the A2 index cell.
$totalWs.Range("B1").Copy()
$q4Ws.Range("B1:H1").PasteSpecial($xlPasteFormats)
$totalWs.Range("A2").Copy()
$q4Ws.Range("A2").PasteSpecial($xlPasteFormats)

# Setting text via a leading apostrophe stamps a quote-prefixed style on the
# cell; strip that back off by re-pasting the plain (unstyled) format from a
# cell that was never touched, so these land back on the default style.
$q4Ws.Range("C2").Copy()
$q4Ws.Range("B2").PasteSpecial($xlPasteFormats)
$q4Ws.Range("D2").PasteSpecial($xlPasteFormats)
$q4Ws.Range("E2").PasteSpecial($xlPasteFormats)
$q4Ws.Range("F2").PasteSpecial($xlPasteFormats)
$q4Ws.Range("G2").PasteSpecial($xlPasteFormats)

# Match page margins used elsewhere in this workbook (0.75/0.75/1/1/0.5/0.5 in).
$q4Ws.PageSetup.LeftMargin = 54
$q4Ws.PageSetup.RightMargin = 54
$q4Ws.PageSetup.TopMargin = 72
$q4Ws.PageSetup.BottomMargin = 72
$q4Ws.PageSetup.HeaderMargin = 36
$q4Ws.PageSetup.FooterMargin = 36

# Adding the sheet makes it active; restore the original active tab
# ("2022-Q2", now the 3rd sheet) via a freshly-resolved reference.
$wb.Worksheets.Item(3).Activate()
